# Actualización automática 2025-07-28 15:10:09
# New advisor "CARAVEDO PAZMIÑO  JAHAIRA PAMELA" is inserted (in alphabetical
# order) between "AVILA TORRES RAFAEL ALEJANDRO" and "EQUISAB S.A.", i.e. as
# the new row 6, on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets.
# This pushes every advisor below down by one row and grows the trailing
# summary row by one.

$wb = $excel.ActiveWorkbook

$newName = "CARAVEDO PAZMIÑO  JAHAIRA PAMELA"

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R, data rows 2-14, summary row 15)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a fresh row above the old row 6 ("EQUISAB S.A."), shifting rows
# 6:15 down to 7:16 (data + the "N de 13" summary row move intact).
$ws1.Rows("6:6").Insert()

$ws1.Range("A6").Value = "OFICINA-CATAECSA"
$ws1.Range("B6").Value = $newName
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(6, $c).Value = 0
}

# The trailing summary row (now row 16) held static "N de 13" labels that
# don't recompute on their own - refresh them to "N de 14" (14 advisors now).
$counts1 = @{ 3=0; 4=0; 5=1; 6=0; 7=0; 8=0; 9=1; 10=0; 11=0; 12=2; 13=1; 14=0; 15=1; 16=0; 17=0; 18=0 }
foreach ($c in $counts1.Keys) {
    $ws1.Cells.Item(16, $c).Value = "$($counts1[$c]) de 14"
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G, data rows 2-14, summary row 15)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows("6:6").Insert()

$ws2.Range("A6").Value = "OFICINA-CATAECSA"
$ws2.Range("B6").Value = $newName
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(6, $c).Value = 0
}

# The SUM totals row (now row 16) carries pure numeric totals that are
# unaffected by the extra (all-zero) row, so no values need patching there.
